$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")
$ws.Range("M4").Value2 = "FIAT"
$ws.Range("M5").Value2 = "CHEVROLET"
$ws.Range("M6").Value2 = "VOLKSWAGEN"
$ws.Range("M7").Value2 = "FORD"
$ws.Range("M8").Value2 = "RENAULT"
$ws.Range("N4").Value2 = 248
$ws.Range("N5").Value2 = 237
$ws.Range("N6").Value2 = 193
$ws.Range("N7").Value2 = 136
$ws.Range("N8").Value2 = 108
$wb.RefreshAll()
$excel.CalculateFull()
Write-Host "done"
